# Extend the daily-data table with rows for 28/05/2021 .. 28/06/2021
# (date serials 44344 .. 44375), continuing the existing pattern of
# zero-valued rows (aggiornamento fino a 28/06 incluso).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastExistingRow = 269
$firstNewRow = 270
$lastNewRow = 301
$firstSerial = 44344

for ($row = $firstNewRow; $row -le $lastNewRow; $row++) {
    $serial = $firstSerial + ($row - $firstNewRow)

    # Copy the formatting (date style, borders, alignment, number format)
    # from the last existing data row so the new A-column cells match the
    # rest of the column exactly, then set the actual values.
    $ws.Range("A$lastExistingRow").Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}
